# residual block changed to bn-con-bn
# Update the per-row metric values in columns D:G (rows 1-9) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    1 = @{ D = 73.99634758930016;  E = 81.15917896601891;  F = 81.2465081557779;   G = 50.48582899066636 }
    2 = @{ D = 82.64315965187087;  E = 91.79104477611941;  F = 89.0788227794947;   G = 56.96652626645226 }
    3 = @{ D = 82.37175119075714;  E = 82.10526315789474;  F = 88.85717287867;     G = 58.34935759910164 }
    4 = @{ D = 65.99759979096662;  E = 89.58601855256981;  F = 84.44585712280185;  G = 59.50711060001895 }
    5 = @{ D = 83.38054610977025;  E = 93.33333333333333;  F = 92.13180261186042;  G = 61.80863502204452 }
    6 = @{ D = 60.34451962776261;  E = 70;                 F = 69.46200917038148;  G = 30.90858859204694 }
    7 = @{ D = 82.77008847366228;  E = 73.68421052631578;  F = 77.83944508861406;  G = 43.55886154819326 }
    8 = @{ D = 63.65253892605208;  E = 51.58543628872329;  F = 62.15223808968322;  G = 34.63078129070227 }
    9 = @{ D = 70.81057694355933;  E = 97.18812509319491;  F = 86.00471750471749;  G = 58.15677100677099 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    $ws.Range("D$row").Value = $cols.D
    $ws.Range("E$row").Value = $cols.E
    $ws.Range("F$row").Value = $cols.F
    $ws.Range("G$row").Value = $cols.G
}
